$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.054.63'
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').Value = '2.311.87'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '301.86'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.03'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.92%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.508'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.64%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.521'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.26%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.81'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.49%  '
$ws.Range('E11').Value = '  -1.02%  '
$ws.Range('E12').Value = '  -1.12%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '17.81'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.82%  '
$ws.Range('E14').Value = '  +0.00%  '
$ws.Range('D15').Value = '2.673.30'
$ws.Range('E15').Value = '  -0.58%  '
$ws.Range('D16').Value = '2.264.19'
$ws.Range('E16').Value = '  -2.60%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.788'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.26%  '
$ws.Range('D18').Value = '42.974.36'
$ws.Range('E18').Value = '  -0.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.30'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.19'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.48%  '
$ws.Range('D21').Value = '0.0₃0907'
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.10'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.50%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.96'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.16'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.98%  '
$ws.Range('E25').Value = '  -1.27%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.11'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.97%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '169.24'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.73%  '
$ws.Range('E29').Value = '  -2.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.16'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.24%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '33.46'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.94%  '
$ws.Range('E32').Value = '  +6.57%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.20'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.32%  '
$ws.Range('E34').Value = '  -0.12%  '
$ws.Range('E35').Value = '  +7.09%  '
$ws.Range('E36').Value = '  -0.79%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0695'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.35%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.103'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.59%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.82'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.04%  '
$ws.Range('E40').Value = '  -2.55%  '
$ws.Range('E41').Value = '  -0.43%  '
$ws.Range('D42').Value = '1.993.44'
$ws.Range('E42').Value = '  -0.56%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0289'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.96%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.12'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.60%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '17.45'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.52%  '
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.84'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.43%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.02'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -13.37%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '76.07'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +8.38%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '54.82'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.71%  '
$ws.Range('D50').Value = '2.540.72'
$ws.Range('E50').Value = '  +0.76%  '
$ws.Range('E51').Value = '  -0.08%  '
